# Applies the OOXML diff to the TBE/doc/todo.docx "ShareFrame" TODO table.
# Strategy: for each affected cell paragraph, locate it by its (currently)
# unique text via Find, then replace the whole paragraph with freshly
# authored OOXML (preserving original <w:p> attributes / <w:pPr>) that
# matches the after-state of the diff, including <w:proofErr> spell-check
# markers split around the camel-case / foreign words.

$d = $word.ActiveDocument

function Replace-ParaText {
    param($SearchFrom, $FindText, $NewXml)
    $rng = $d.Content.Duplicate
    $rng.Start = $SearchFrom
    $rng.End = $d.Content.End
    $ok = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) { throw "Could not find text: $FindText" }
    $rng.InsertXML($NewXml)
    return $rng.End
}

$cursor = 0

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="007D6018" w:rsidRPr="007D6018" w:rsidRDefault="007D6018" w:rsidP="007D6018"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">ShareFrame: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>JTree</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> vom Server</w:t></w:r></w:p>'
$cursor = Replace-ParaText $cursor "ShareFrame: JTree vom Server" $newXml

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="007D6018" w:rsidRPr="007D6018" w:rsidRDefault="007D6018" w:rsidP="007D6018"><w:proofErr w:type="spellStart"/><w:r><w:t>ShareFrame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Buttons („upload“, „download“, „Cancel“, „Open“)</w:t></w:r></w:p>'
$cursor = Replace-ParaText $cursor "ShareFrame: Buttons („upload“, „download“, „Cancel“, „Open“)" $newXml

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="007D6018" w:rsidRPr="009368CD" w:rsidRDefault="007D6018" w:rsidP="009368CD"><w:proofErr w:type="spellStart"/><w:r><w:t>SettingsFrame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Usereingaben</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>überprüfen</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$cursor = Replace-ParaText $cursor "SettingsFrame: Usereingaben überprüfen" $newXml

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="007D6018" w:rsidRPr="007D6018" w:rsidRDefault="007D6018" w:rsidP="007D6018"><w:r><w:t xml:space="preserve">FTP: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>progressMonitoring</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>für</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ShareFrame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> und </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SettingsFrame</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$cursor = Replace-ParaText $cursor "FTP: progressMonitoring für ShareFrame und SettingsFrame" $newXml

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="007D6018" w:rsidRPr="007D6018" w:rsidRDefault="007D6018" w:rsidP="007D6018"><w:proofErr w:type="spellStart"/><w:r><w:t>SettingsFrame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Sportarten</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> / </w:t></w:r><w:r><w:t xml:space="preserve">Languages </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>downloaden</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$cursor = Replace-ParaText $cursor "SettingsFrame: Sportarten / Languages downloaden" $newXml

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="007D6018" w:rsidRPr="007D6018" w:rsidRDefault="007D6018" w:rsidP="007D6018"><w:r><w:t xml:space="preserve">FTP: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ExceptionHandling</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$cursor = Replace-ParaText $cursor "FTP: ExceptionHandling" $newXml

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="007D6018" w:rsidRPr="007D6018" w:rsidRDefault="007D6018" w:rsidP="007D6018"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Refreshing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> aller </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Gui</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>-Elemente optimieren</w:t></w:r></w:p>'
$cursor = Replace-ParaText $cursor "Refreshing aller Gui-Elemente optimieren" $newXml

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="007D6018" w:rsidRPr="007D6018" w:rsidRDefault="007D6018" w:rsidP="007D6018"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">Shapes auf Board </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>resizen</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$cursor = Replace-ParaText $cursor "Shapes auf Board resizen" $newXml

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00123D9E" w:rsidRDefault="00363785" w:rsidP="007D6018"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>CurvedBezierTool</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$cursor = Replace-ParaText $cursor "CurvedBezierTool" $newXml

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00EF245F" w:rsidRDefault="00445AE4" w:rsidP="007D6018"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Menukonzept</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="007F147C"><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Menus deaktivieren</w:t></w:r></w:p>'
$cursor = Replace-ParaText $cursor "Menukonzept: Menus deaktivieren" $newXml

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00EF245F" w:rsidRDefault="00445AE4" w:rsidP="007D6018"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>About</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$cursor = Replace-ParaText $cursor "About" $newXml

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="003653AD" w:rsidRDefault="003653AD" w:rsidP="007D6018"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Rotatebild</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$cursor = Replace-ParaText $cursor "Rotatebild" $newXml

# Empty Wingdings-checkbox cell right after "Attributliste nicht nur Titel
# anzeigen" gains a checkmark run (table cell 24, column 2).
$symXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="003653AD" w:rsidRPr="007E1874" w:rsidRDefault="003653AD" w:rsidP="007D6018"><w:pPr><w:cnfStyle w:val="000000100000"/><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="de-CH"/></w:rPr><w:sym w:font="Wingdings" w:char="F0FC"/></w:r></w:p>'
$symCell = $d.Tables.Item(1).Cell(24, 2)
$symRange = $symCell.Range
$symRange.InsertXML($symXml)
$cursor = $symCell.Range.End

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009804C1" w:rsidRDefault="00BF7429" w:rsidP="007D6018"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>ExceptionHandling</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$cursor = Replace-ParaText $cursor "ExceptionHandling" $newXml

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00550D71" w:rsidRDefault="00550D71" w:rsidP="007D6018"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Recently</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>Opened</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> Files aktualisieren</w:t></w:r></w:p>'
$cursor = Replace-ParaText $cursor "Recently Opened Files aktualisieren" $newXml

Write-Output "done"
